# Updates the cryptocurrency price/volume table to reflect the latest
# symbol-list snapshot (GitHub Actions scheduled refresh).
#
# For every affected row, columns D (Price) and E (Volume(1h)) are updated.
# Both columns are stored as text in the sheet (e.g. "330.90", "1.08%"),
# so we explicitly force the cell's number format to Text ("@") before
# writing the value. This prevents Excel from re-interpreting the text as
# a number (which would, for example, turn "330.90" into 330.9 and drop
# the trailing zero, or turn "1.08%" into a numeric percentage).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "330.90"; E = "1.08%" }
    @{ Row = 3; D = "40.13"; E = "0.89%" }
    @{ Row = 4; D = "5.716"; E = "0.43%" }
    @{ Row = 5; D = "0.08111"; E = "0.62%" }
    @{ Row = 6; D = "8.655"; E = "-0.29%" }
    @{ Row = 7; D = "4.485"; E = "-1.80%" }
    @{ Row = 8; D = "1.958"; E = "0.44%" }
    @{ Row = 9; D = "2.953"; E = "-0.85%" }
    @{ Row = 10; D = "0.9287"; E = "-1.64%" }
    @{ Row = 11; D = "0.1260"; E = "-1.29%" }
    @{ Row = 12; D = "0.1958"; E = "-1.40%" }
    @{ Row = 13; D = "8.725"; E = "14.39%" }
    @{ Row = 14; D = "0.09173"; E = "-0.72%" }
    @{ Row = 15; D = "0.03571"; E = "0.75%" }
    @{ Row = 16; D = "0.1051"; E = "9.46%" }
    @{ Row = 17; D = "0.001295"; E = "-0.83%" }
    @{ Row = 18; D = "0.006299"; E = "2.37%" }
    @{ Row = 19; D = "3.376"; E = "-0.04%" }
    @{ Row = 20; D = "0.3448"; E = "-1.42%" }
    @{ Row = 21; D = "0.1368"; E = "-3.10%" }
    @{ Row = 22; D = "0.2605"; E = "3.58%" }
    @{ Row = 23; D = "0.04411"; E = "0.22%" }
    @{ Row = 24; D = "0.001249"; E = "-0.46%" }
    @{ Row = 25; D = "0.004467"; E = "3.36%" }
    @{ Row = 26; D = "0.0001235"; E = "3.49%" }
    @{ Row = 39; D = "0.02724"; E = "8.26%" }
    @{ Row = 40; D = "0.05551"; E = "6.72%" }
    @{ Row = 41; D = "0.007515"; E = "2.94%" }
    @{ Row = 42; D = "0.009813"; E = "8.06%" }
    @{ Row = 43; D = "0.1420"; E = "-0.24%" }
    @{ Row = 44; D = "0.002099"; E = "-4.41%" }
    @{ Row = 45; D = "0.01185"; E = "18.09%" }
    @{ Row = 46; D = "0.00006786"; E = "0.59%" }
    @{ Row = 47; D = "0.00000000747"; E = "-0.57%" }
    @{ Row = 48; D = "0.003060"; E = "6.35%" }
    @{ Row = 49; D = "0.002271"; E = "25.90%" }
    @{ Row = 50; D = "0.00002092"; E = "-0.57%" }
    @{ Row = 51; D = "0.0001992"; E = "-0.57%" }
)

foreach ($u in $updates) {
    $dCell = $ws.Cells.Item($u.Row, 4)   # Column D = Price
    $eCell = $ws.Cells.Item($u.Row, 5)   # Column E = Volume(1h)

    $dCell.NumberFormat = "@"
    $dCell.Value = $u.D

    $eCell.NumberFormat = "@"
    $eCell.Value = $u.E
}
